$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (also updates the tab name / workbook.xml <sheet name=.../>)
$ws.Name = "Through 2022-07-11"

# Update the label for the July row
$ws.Range("A8").Value = "July (through 07-11)"

# Update July (row 8) values
$ws.Range("B8").Value = 15
$ws.Range("C8").Value = 20
$ws.Range("D8").Value = 19
$ws.Range("E8").Value = 28
$ws.Range("F8").Value = 16
$ws.Range("G8").Value = 39
$ws.Range("H8").Value = 54
$ws.Range("I8").Value = 63

# Update Total (row 9) values
$ws.Range("B9").Value = 140
$ws.Range("C9").Value = 268
$ws.Range("D9").Value = 409
$ws.Range("E9").Value = 381
$ws.Range("F9").Value = 267
$ws.Range("G9").Value = 511
$ws.Range("H9").Value = 814
$ws.Range("I9").Value = 868

$wb.Save()
